# Adds the "sigma_ratio" data column (C) and a filtered "Name" helper
# column (Q) to Sheet1, per commit "added sigma_ratio to excel".
#
# Column C (rows 2-36): numeric sigma_ratio value per source, or "-" for
# the handful of rows whose State (B) marks them not yet fit (their
# D/M columns already show the "to do" / reject markers).
#
# Column Q (rows 2-31): the subset of Name (column A) values for which a
# sigma_ratio was actually computed (i.e. excludes the "-" rows), listed
# in original order - effectively a "good" filter pasted as values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column C: sigma_ratio -------------------------------------------------
$sigmaRatioData = @(
    @{Row=2; Val=78.654844507959666},
    @{Row=3; Val=65.698367298851593},
    @{Row=4; Val=53.953585879898611},
    @{Row=5; Val=59.80726881378898},
    @{Row=6; Val=66.486284026074728},
    @{Row=7; Val=82.504263953251424},
    @{Row=8; Val=64.840124127438841},
    @{Row=9; Val="-"},
    @{Row=10; Val="-"},
    @{Row=11; Val=69.873372668561842},
    @{Row=12; Val=36.268930152835487},
    @{Row=13; Val=90.542280136417446},
    @{Row=14; Val=87.693954714566871},
    @{Row=15; Val=81.947799671078101},
    @{Row=16; Val=51.447760576429062},
    @{Row=17; Val=45.471096843211953},
    @{Row=18; Val=68.342485013728535},
    @{Row=19; Val=63.807178135337203},
    @{Row=20; Val=31.692075915159261},
    @{Row=21; Val=89.123003413701312},
    @{Row=22; Val=80.271706068618229},
    @{Row=23; Val="-"},
    @{Row=24; Val=89.745617229461871},
    @{Row=25; Val=82.611830801977987},
    @{Row=26; Val=41.381685153983703},
    @{Row=27; Val=76.395308522288346},
    @{Row=28; Val="-"},
    @{Row=29; Val=70.553748491929099},
    @{Row=30; Val=71.860842600437707},
    @{Row=31; Val="-"},
    @{Row=32; Val=80.096582091956719},
    @{Row=33; Val=82.389461370168149},
    @{Row=34; Val=64.873489847899037},
    @{Row=35; Val=70.815824382423372},
    @{Row=36; Val=61.312945662889689}
)

foreach ($item in $sigmaRatioData) {
    $ws.Cells.Item($item.Row, 3).Value = $item.Val
}

# ---- Column Q: filtered Name list ------------------------------------------
$filteredNameData = @(
    @{Row=2; Val="ESO079-003_GROUP_factor2.5_pixscale0.6"},
    @{Row=3; Val="ESO298-016_factor2.5_pixscale0.6"},
    @{Row=4; Val="ESO356-012_factor2.5_pixscale0.6"},
    @{Row=5; Val="IC0120_factor2.5_pixscale0.6"},
    @{Row=6; Val="IC1124_factor2.5_pixscale0.6"},
    @{Row=7; Val="IC1657_GROUP_factor10.0_pixscale0.6"},
    @{Row=8; Val="NGC0804_factor3.0_pixscale0.6"},
    @{Row=9; Val="NGC1162_factor2.5_pixscale0.6"},
    @{Row=10; Val="NGC2769_factor2.5_pixscale0.6"},
    @{Row=11; Val="NGC3041_factor2.5_pixscale0.6"},
    @{Row=12; Val="NGC4256_GROUP_factor22.0_pixscale0.6"},
    @{Row=13; Val="NGC4388_GROUP_factor16.0_pixscale0.6"},
    @{Row=14; Val="NGC4632_factor2.5_pixscale0.6"},
    @{Row=15; Val="NGC4686_factor2.5_pixscale0.6"},
    @{Row=16; Val="NGC5055_GROUP_factor75.0_pixscale0.6"},
    @{Row=17; Val="NGC5263_factor2.5_pixscale0.6"},
    @{Row=18; Val="NGC5387_factor2.5_pixscale0.6"},
    @{Row=19; Val="NGC5513_GROUP_factor7.0_pixscale0.6"},
    @{Row=20; Val="NGC5907_factor2.5_pixscale0.6"},
    @{Row=21; Val="PGC006791_factor2.5_pixscale0.6"},
    @{Row=22; Val="PGC021008_factor2.5_pixscale0.6"},
    @{Row=23; Val="PGC039258_factor2.5_pixscale0.6"},
    @{Row=24; Val="PGC1001085_factor3.5_pixscale0.6"},
    @{Row=25; Val="PGC3092153_factor2.5_pixscale0.6"},
    @{Row=26; Val="PGC430221_factor4.0_pixscale0.6"},
    @{Row=27; Val="UGC01245_factor4.5_pixscale0.6"},
    @{Row=28; Val="UGC01424_factor3.0_pixscale0.6"},
    @{Row=29; Val="UGC01970_factor3.5_pixscale0.6"},
    @{Row=30; Val="UGC08717_factor3.0_pixscale0.6"},
    @{Row=31; Val="UGC09239_factor2.5_pixscale0.6"}
)

foreach ($item in $filteredNameData) {
    $ws.Cells.Item($item.Row, 17).Value = $item.Val
}

# ---- Column widths for the two new columns (Q, R) --------------------------
$ws.Columns.Item(17).ColumnWidth = 36.166666666666664
$ws.Columns.Item(18).ColumnWidth = 11.330729166666666

# ---- Restore the selected cell to match the saved view ---------------------
$ws.Range("I9").Select()
